# Add a new record (row 4) to the "baseDatosPersonas" table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Felipene Obando"
$ws.Range("B4").Value = "felipeobando2001@gmail.com"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "San José,Alajuela,Heredia,Puntarenas,Guanacaste,Cartago,Limón"
